# Commit: "changing document, table attributes to lowerCamelCase"
#
# The ObjTables header rows encode attribute-like tokens (ObjTablesVersion,
# Type, Id) inside plain text cell values. This change renames those tokens
# to lowerCamelCase (objTablesVersion, type, id) across the three sheets.

$wb = $excel.ActiveWorkbook

$tocSheet  = $wb.Worksheets.Item("!!_Table of contents")
$dataSheet = $wb.Worksheets.Item("!!Data repo metadata")
$modelSheet = $wb.Worksheets.Item("!!Model1s")

$tocSheet.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$tocSheet.Range("A2").Value = "!!ObjTables type='TableOfContents'"

$dataSheet.Range("A1").Value = "!!ObjTables type='Data' id='DataRepoMetadata'"

$modelSheet.Range("A1").Value = "!!ObjTables type='Data' id='Model1'"
